$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "eADWE152"
$ws.Range("B2").Value = 23111045
$ws.Range("C2").Value = "swlkjgw87"
$ws.Range("D2").Value = "ChN2%#b6"
$ws.Range("F2").Value = "BzEgePuf"
$ws.Range("G2").Value = "SsIu"
